$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the old row 537 (weekly update: new market-day
# records for Packham's Triumph, "Primera" and "Segunda" quality grades).
# This shifts all former rows 537:592 down to 539:594.
$ws.Rows("537:538").Insert()

# New row 537: Packham's Triumph, Primera
$ws.Range("A537").Value = 4
$ws.Range("B537").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C537").Value = "Los Lagos"
$ws.Range("D537").Value = "2023-09-25"
$ws.Range("E537").Value = 10
$ws.Range("F537").Value = "Fruta"
$ws.Range("G537").Value = 100104
$ws.Range("H537").Value = "Frutos de pepita"
$ws.Range("I537").Value = 100104005
$ws.Range("J537").Value = "Pera"
$ws.Range("K537").Value = "Packham's Triumph"
$ws.Range("L537").Value = "Primera"
$ws.Range("M537").Value = 200
$ws.Range("N537").Value = 20000
$ws.Range("O537").Value = 20000
$ws.Range("P537").Value = 20000
$ws.Range("Q537").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R537").Value = "Región de O'Higgins"
$ws.Range("S537").Value = 1333
$ws.Range("T537").Value = 15

# New row 538: Packham's Triumph, Segunda
$ws.Range("A538").Value = 4
$ws.Range("B538").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C538").Value = "Los Lagos"
$ws.Range("D538").Value = "2023-09-25"
$ws.Range("E538").Value = 10
$ws.Range("F538").Value = "Fruta"
$ws.Range("G538").Value = 100104
$ws.Range("H538").Value = "Frutos de pepita"
$ws.Range("I538").Value = 100104005
$ws.Range("J538").Value = "Pera"
$ws.Range("K538").Value = "Packham's Triumph"
$ws.Range("L538").Value = "Segunda"
$ws.Range("M538").Value = 200
$ws.Range("N538").Value = 16000
$ws.Range("O538").Value = 16000
$ws.Range("P538").Value = 16000
$ws.Range("Q538").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R538").Value = "Región de O'Higgins"
$ws.Range("S538").Value = 1067
$ws.Range("T538").Value = 15
